$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 42
$ws.Range("F2").Value = 24
$ws.Range("H2").Value = 24

$ws.Range("E15").Value = 146
$ws.Range("F15").Value = 76
$ws.Range("H15").Value = 76

$ws.Range("E18").Value = 102

$ws.Range("E36").Value = 90

$ws.Range("E45").Value = 23
$ws.Range("F45").Value = 12
$ws.Range("H45").Value = 12

$ws.Range("F48").Value = 16
$ws.Range("H48").Value = 16

$ws.Range("E49").Value = 57

$ws.Range("E64").Value = 32

$ws.Range("E66").Value = 31

$ws.Range("E70").Value = 39

$ws.Range("E73").Value = 29

$ws.Range("E76").Value = 45

$ws.Range("E79").Value = 31
